$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: remove the "assists" entry that used to live in L8 ---
$ws.Range("L8").ClearContents()

# --- Row 12: add a new "got_refrafgt" label in L12 (new shared string) ---
$ws.Range("L12").Value = "got_refrafgt"

# --- New vertical transcription of the PlayerRound header row (row 8)
#     into column C, rows 34-54. C34 carries the same bold style as B8. ---
$ws.Range("C34").Value = "PlayerRound"
$ws.Range("C34").Font.Bold = $true

$ws.Range("C35").Value = "id[serial]"
$ws.Range("C36").Value = "Player.ID"
$ws.Range("C37").Value = "Round.ID"
$ws.Range("C38").Value = "Operator"
$ws.Range("C39").Value = "Team Index"
$ws.Range("C40").Value = "spawn"
$ws.Range("C41").Value = "kills"
$ws.Range("C42").Value = "death"
$ws.Range("C43").Value = "assists"
$ws.Range("C44").Value = "headshots"
$ws.Range("C45").Value = "plant"
$ws.Range("C46").Value = "defuse"
$ws.Range("C47").Value = "kostpoint"
$ws.Range("C48").Value = "1vX /selber rechnen"
$ws.Range("C49").Value = "ok"
$ws.Range("C50").Value = "od"
$ws.Range("C51").Value = "win"
$ws.Range("C52").Value = "ATK"
$ws.Range("C53").Value = "refrags"
$ws.Range("C54").Value = "got refragt"

# --- Page setup (portrait, paper size 9 = A4) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection / view: active cell moves to L16, no frozen top-left cell ---
$ws.Range("L16").Select()
